# Chapter 7 Answer Key — homework edit
#  1. Exercise 13: append the rest of the Groucho Marx joke.
#  2. Exercise 14: insert a new part (c) that asks for two bracket-notation
#     diagrams (garden-path vs. correct reading), reletter the old part (c)
#     "Model response" to part (d).
#  3. Remove Exercise 15 (and its "Model response" sub-heading) entirely,
#     re-using its old text as the new part (d) answer.

$d = $word.ActiveDocument

# Helper: scan the (stable) Paragraphs collection by plain-text match,
# starting at paragraph index $startAt. Paragraph ranges include the
# trailing paragraph mark, so trim before comparing. This is far more
# reliable in this host than Find-narrowed ranges / Paragraphs.First.Index.
function Find-ParaIndex {
    param($startAt, $text)
    for ($i = $startAt; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd()
        if ($t -eq $text) { return $i }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) Exercise 13 — extend the Groucho Marx quotation.
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "I shot an elephant in my pajamas.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "I shot an elephant in my pajamas. How he got in my pajamas, I will never know.",
    2)

# ---------------------------------------------------------------------
# 2) Locate Exercise 14's "c) Model response:" heading (there are two
#    "c) Model response:" paragraphs in the document — one belongs to
#    Exercise 13 and must stay untouched — so anchor the search on the
#    unique "b) Correct reading:" paragraph that immediately precedes
#    Exercise 14's copy.
# ---------------------------------------------------------------------
$bCorrectIndex = Find-ParaIndex 1 "b) Correct reading:"
$cModelIndex = Find-ParaIndex ($bCorrectIndex + 1) "c) Model response:"

# 2a) Relabel "c) Model response:" -> "c) Bracket notation for each reading:"
$cModelPara = $d.Paragraphs.Item($cModelIndex)
$cModelPara.Range.Text = "c) Bracket notation for each reading:"

# ---------------------------------------------------------------------
# 3) Insert the "Diagram 1" paragraph right after the relabeled part (c)
#    heading, ahead of the garden-path explanation paragraph.
# ---------------------------------------------------------------------
$gardenIndex = $cModelIndex + 1
$gardenPara = $d.Paragraphs.Item($gardenIndex)
$gardenPara.Range.InsertParagraphBefore()

$diagram1Index = $gardenIndex
$diagram1Para = $d.Paragraphs.Item($diagram1Index)
$diagram1Range = $diagram1Para.Range
$diagram1Range.Text = "Diagram 1 " + [char]0x2014 + " Garden-path (incorrect) reading: "
$diagram1BoldRange = $d.Range($diagram1Range.Start, $diagram1Range.End - 1)
$diagram1BoldRange.Font.Bold = $true
$diagram1BoldRange.Font.Size = 12

$diagram1Para = $d.Paragraphs.Item($diagram1Index)
$insertAt = $diagram1Para.Range.End - 1
$insertPoint = $d.Range($insertAt, $insertAt)
$bracket1 = "[S [NP [DET The] [N horse]] [VP [V raced] [PP [PREP past] [NP [DET the] [N barn]]]]] + fell ???"
$insertPoint.InsertAfter($bracket1)
$bracket1Range = $d.Range($insertAt, $insertAt + $bracket1.Length)
$bracket1Range.Font.Name = "Consolas"
$bracket1Range.Font.Size = 11
$bracket1Range.Font.Bold = $false

# ---------------------------------------------------------------------
# 4) Replace the garden-path explanation paragraph's text (it shifted
#    down by one position after the Diagram 1 insertion above).
# ---------------------------------------------------------------------
$gardenIndex = $diagram1Index + 1
$gardenPara = $d.Paragraphs.Item($gardenIndex)
$gardenPara.Range.Text = "In the garden-path reading, " + [char]34 + "raced" + [char]34 + " is parsed as the main verb with " + [char]34 + "past the barn" + [char]34 + " as a PP inside the VP. This leaves " + [char]34 + "fell" + [char]34 + " with no grammatical role, which is why the sentence seems to break."

# ---------------------------------------------------------------------
# 5) Remove "Exercise 15." heading paragraph outright (its old spot is
#    replaced by two brand-new paragraphs below, so this just drops it
#    and its distinctive "before=120/after=60" spacing).
# ---------------------------------------------------------------------
$ex15Index = $gardenIndex + 1
$ex15Para = $d.Paragraphs.Item($ex15Index)
$ex15Para.Range.Delete()

# ---------------------------------------------------------------------
# 6) Insert two fresh paragraphs after the garden-path explanation:
#      - "Diagram 2 -- Correct reading: " + bracket notation
#      - "In the correct reading, ..." explanation
#    Inserting both via the (still valid) garden-path paragraph range
#    means they inherit its clean "ind left=1008" formatting exactly.
# ---------------------------------------------------------------------
$gardenPara = $d.Paragraphs.Item($gardenIndex)
$gardenPara.Range.InsertParagraphAfter()
$gardenPara.Range.InsertParagraphAfter()

$diagram2Index = $gardenIndex + 1
$diagram2Para = $d.Paragraphs.Item($diagram2Index)
$diagram2Range = $diagram2Para.Range
$diagram2Range.Text = "Diagram 2 " + [char]0x2014 + " Correct reading: "
$diagram2BoldRange = $d.Range($diagram2Range.Start, $diagram2Range.End - 1)
$diagram2BoldRange.Font.Bold = $true
$diagram2BoldRange.Font.Size = 12

$diagram2Para = $d.Paragraphs.Item($diagram2Index)
$insertAt2 = $diagram2Para.Range.End - 1
$insertPoint2 = $d.Range($insertAt2, $insertAt2)
$bracket2 = "[S [NP [DET The] [N horse] [VP [V raced] [PP [PREP past] [NP [DET the] [N barn]]]]] [VP [V fell]]]"
$insertPoint2.InsertAfter($bracket2)
$bracket2Range = $d.Range($insertAt2, $insertAt2 + $bracket2.Length)
$bracket2Range.Font.Name = "Consolas"
$bracket2Range.Font.Size = 11
$bracket2Range.Font.Bold = $false

$correctReadingIndex = $diagram2Index + 1
$correctReadingPara = $d.Paragraphs.Item($correctReadingIndex)
$correctReadingPara.Range.Text = "In the correct reading, " + [char]34 + "raced past the barn" + [char]34 + " is a reduced relative clause inside the subject NP (modifying " + [char]34 + "horse" + [char]34 + "), and " + [char]34 + "fell" + [char]34 + " is the main verb of the sentence."

# ---------------------------------------------------------------------
# 7) Reletter the trailing "Model response:" heading to "d) Model
#    response:" and replace Exercise 15's old model-response text with
#    the garden-path essay text that used to live in part (c).
# ---------------------------------------------------------------------
$dModelIndex = $correctReadingIndex + 1
$dModelPara = $d.Paragraphs.Item($dModelIndex)
$dModelPara.Range.Text = "d) Model response:"

$finalIndex = $dModelIndex + 1
$finalPara = $d.Paragraphs.Item($finalIndex)
$finalPara.Range.Text = "Garden-path sentences cause confusion because our brains process language incrementally " + [char]0x2014 + " we build structural interpretations word by word as we read. When we encounter " + [char]34 + "The horse raced," + [char]34 + " the simplest analysis is that " + [char]34 + "raced" + [char]34 + " is the main verb, and we commit to that structure. When " + [char]34 + "fell" + [char]34 + " appears, it forces us to revise: " + [char]34 + "raced" + [char]34 + " was actually part of a reduced relative clause, not the main verb. This revision is cognitively costly, which is why the sentence feels confusing. Garden-path sentences demonstrate that sentence comprehension is not just about knowing the words " + [char]0x2014 + " it requires actively building and sometimes revising hierarchical structure in real time."

Write-Output "Done. Total paragraphs: $($d.Paragraphs.Count)"
